# Added Members Collect method
# Fill in the placeholder "???" pricing statuses on the PricingStatus sheet
# and make that sheet the active tab/selection.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PricingStatus")

$ws.Range("B3").Value = "Hourly"
$ws.Range("B4").Value = "PerCapita"
$ws.Range("B5").Value = "ProBono"

$ws.Activate()
$ws.Range("C5").Select()
